$d = $word.ActiveDocument

$replacements = @(
    @{old="624×6="; new="833×3="},
    @{old="352×5="; new="470×2="},
    @{old="314×6="; new="941×8="},
    @{old="622×7="; new="978×2="},
    @{old="410×8="; new="179×8="},
    @{old="837×9="; new="101×6="},
    @{old="707×7="; new="828×5="},
    @{old="150×5="; new="514×9="},
    @{old="852×4="; new="647×7="},
    @{old="871×8="; new="924×3="},
    @{old="616×3="; new="870×7="},
    @{old="142×4="; new="626×5="},
    @{old="379×2="; new="179×4="},
    @{old="858×7="; new="677×5="},
    @{old="395×4="; new="949×8="},
    @{old="773×7="; new="975×8="},
    @{old="828×9="; new="355×7="},
    @{old="137×4="; new="770×2="},
    @{old="187×3="; new="249×3="},
    @{old="259×6="; new="567×9="},
    @{old="700×3="; new="822×3="},
    @{old="647×5="; new="421×9="},
    @{old="300×2="; new="627×3="},
    @{old="672×4="; new="797×8="},
    @{old="152×2="; new="836×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
